$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data (rows 2-15) after fixing/reordering the geopoints and adding a new
# labeled point ("little field fountain") that duplicates point #1's coordinates.
$data = @(
    @("little field fountain", 30.2837284764915,   -97.739589214324894),
    @(1,                       30.2837284764915,   -97.739589214324894),
    @(2,                       30.285441285189599, -97.735724821686702),
    @(3,                       30.283015086901401, -97.737016975879598),
    @(4,                       30.286693156878201, -97.740957140922504),
    @(5,                       30.2862808863749,   -97.736974060535402),
    @(6,                       30.289081049140002, -97.740678191184898),
    @(7,                       30.283865131697599, -97.7382239699363),
    @(8,                       30.2867811694578,   -97.740332186222005),
    @(9,                       30.288407071857801, -97.736746072769094),
    @(10,                      30.283450533112099, -97.738524377345996),
    @(11,                      30.2840434780362,   -97.736920416355105),
    @(12,                      30.287693716288,    -97.739219069480896),
    @(13,                      30.286125705455699, -97.737349569797502)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# One fewer data row than before (15 -> 14), so delete the now-unused last row.
$ws.Rows.Item(16).Delete()

$ws.Range("A19").Select()
